# Append a new row (row 9) to Sheet1, mirroring the existing "احمد / الجزائري / الرحلة 1 / C3 / NRC"
# entries already present in rows 6-8, with a new quantity and timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Column A (ملاحظات / notes) - left blank like most existing rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = ""
$ws.Cells.Item($row, 1).Style = "Normal"

# Column B (المرافق / accompanying person)
$ws.Cells.Item($row, 2).Value = "احمد"

# Column C (الكمية / quantity) - stored as text, matching the rest of the column.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "6"
$ws.Cells.Item($row, 3).Style = "Normal"

# Column D (المخيم / camp)
$ws.Cells.Item($row, 4).Value = "الجزائري"

# Column E (نوع المسافة / trip type)
$ws.Cells.Item($row, 5).Value = "الرحلة 1"

# Column F (المركبة / vehicle)
$ws.Cells.Item($row, 6).Value = "C3"

# Column G (المؤسسة / organization)
$ws.Cells.Item($row, 7).Value = "NRC"

# Column H (الوقت / timestamp)
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٤١:٠٤ م"
